# "cambio de ambiente de pruebas, ejercicio parte uno finalizado"
#
# 1) Hoja1 becomes the active/selected sheet (workbook no longer remembers
#    "Credenciales" as the active tab); selection on Hoja1 moves to B1 and
#    its (previously stray) B1 value is cleared out.
# 2) Credenciales gets a fresh set of demo credentials (2 rows instead of a
#    single 4-column row) and a print setup, with the selection left on C2.

$wb = $excel.ActiveWorkbook

$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja2 = $wb.Worksheets.Item("Hoja2")
$wsCred  = $wb.Worksheets.Item("Credenciales")

# --- Hoja1: clear the stray B1 cell, widen column B a bit, active selection B1 ---
$wsHoja1.Range("B1").ClearContents()
$wsHoja1.Columns.Item(2).ColumnWidth = 23.6

# --- Credenciales: replace the old single-row credentials with the new
#     two-row demo user / demo admin set, then clear the now-unused C/D ---
$wsCred.Range("A1").Value = "user@phptravels.com"
$wsCred.Range("B1").Value = "demouser"
$wsCred.Range("B2").Value = "demoadmin"
$wsCred.Range("A2").Value = "admin@phptravels.com"
$wsCred.Range("C1").ClearContents()
$wsCred.Range("D1").ClearContents()
$wsCred.Columns.Item(1).ColumnWidth = 22.6
$wsCred.Columns.Item(3).ColumnWidth = 15.7
$wsCred.PageSetup.Orientation = 1

# --- leave Credenciales' own remembered selection on C2 ---
$wsCred.Range("C2").Select()

# --- activate Hoja1 last so it becomes the workbook's active tab/selection ---
$wsHoja1.Activate()
$wsHoja1.Range("B1").Select()
